# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Caqui" (Fuyu, Primera) at row 6,
# pushing the existing data down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (shifts rows 6..63 down to 7..64)
$ws.Rows.Item(6).Insert(4)

# Populate the newly inserted row 6 with the new weekly record
$ws.Cells.Item(6, 1).Value  = 9
$ws.Cells.Item(6, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(6, 3).Value  = "Metropolitana"
$ws.Cells.Item(6, 4).Value  = 45069
$ws.Cells.Item(6, 5).Value  = 13
$ws.Cells.Item(6, 6).Value  = "Fruta"
$ws.Cells.Item(6, 7).Value  = 100107
$ws.Cells.Item(6, 8).Value  = "Otros"
$ws.Cells.Item(6, 9).Value  = 100107001
$ws.Cells.Item(6, 10).Value = "Caqui"
$ws.Cells.Item(6, 11).Value = "Fuyu"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 470
$ws.Cells.Item(6, 14).Value = 13000
$ws.Cells.Item(6, 15).Value = 13500
$ws.Cells.Item(6, 16).Value = 13234
$ws.Cells.Item(6, 17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(6, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(6, 19).Value = 827
$ws.Cells.Item(6, 20).Value = 16
